# Update the cryptos price/volume table to reflect the latest scrape values.
# Source commit: "Updated cryptos list on Mon Nov 20 05:47:49 UTC 2023 with GitHub Actions"
#
# Cells in column D hold free-text "price" strings (dotted thousands, e.g.
# "37.111.47") that must stay text even when a refreshed value happens to
# look like a plain decimal number (e.g. "245.89"). Setting NumberFormat to
# "@" (Text) before assigning the Value keeps Excel from auto-coercing those
# into numeric cells, and resetting Style back to Normal afterwards avoids
# leaving a stray custom number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.111.47'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '1.998.19'
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  +0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '245.89'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.70%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.625'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.72%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '59.75'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +1.82%  '
$ws.Range('E8').Value = '  +0.10%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.383'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +2.62%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0804'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +2.09%  '
$ws.Range('E11').Value = '  +1.27%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '15.02'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +5.41%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '22.46'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +5.54%  '
$ws.Range('D14').Value = '2.295.18'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.843'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.42%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.42'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +2.42%  '
$ws.Range('D17').Value = '2.002.24'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '37.089.68'
$ws.Range('E18').Value = '  +1.68%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '70.17'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('D20').Value = '0.0₃0862'
$ws.Range('E20').Value = '  +1.56%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.16'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.01%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '230.12'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  +0.13%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.47'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.58%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.35'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.28%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.40'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +2.54%  '
$ws.Range('E27').Value = '  +0.63%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '163.31'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +1.64%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '19.59'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('E30').Value = '  +11.59%  '
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('E32').Value = '  +1.44%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0653'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +6.62%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.48'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.69%  '
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  +2.56%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.29'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -6.37%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '5.37'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.96%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.0980'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.06%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.92'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('E42').Value = '  +1.99%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.18'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.53%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '16.65'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('D46').Value = '1.366.01'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.04'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.73%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.42'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +4.09%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.04'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +11.56%  '
$ws.Range('E50').Value = '  +0.04%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '46.11'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +4.69%  '
